# Weekly update: insert this week's two new price rows (Primera / Segunda)
# for "Sandia" right after the existing row 3, pushing the rest of the
# historical rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 4 (rows 4-30 shift down to 6-32).
$ws.Range("A4:A5").EntireRow.Insert()

# New row 4: "Primera" quality entry dated 2022-02-10
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = "2022-02-10"
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = 100112028
$ws.Range("G4").Value = "Sandia"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 1300
$ws.Range("K4").Value = 350
$ws.Range("L4").Value = 380
$ws.Range("M4").Value = 365
$ws.Range("N4").Value = "$/kilo (volumen en unidades)"
$ws.Range("O4").Value = "Perú"
$ws.Range("P4").Value = 365
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = "Hortaliza"

# New row 5: "Segunda" quality entry, same date
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C5").Value = "Arica y Parinacota"
$ws.Range("D5").Value = "2022-02-10"
$ws.Range("E5").Value = 15
$ws.Range("F5").Value = 100112028
$ws.Range("G5").Value = "Sandia"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 900
$ws.Range("K5").Value = 300
$ws.Range("L5").Value = 330
$ws.Range("M5").Value = 315
$ws.Range("N5").Value = "$/kilo (volumen en unidades)"
$ws.Range("O5").Value = "Perú"
$ws.Range("P5").Value = 315
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = "Hortaliza"
